$d = $word.ActiveDocument
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>Advanced Calculator App - Documentation</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>1. Introduction</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The Advanced Calculator App is a modular Python application designed to perform various arithmetic operations. The project will be built incrementally, starting with basic operations and expanding to advanced features such as scientific functions and a graphical user interface (GUI) using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tkinter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>2. Objectives</w:t></w:r></w:p><w:p><w:r><w:t>• Develop a modular calculator using Python functions.</w:t></w:r></w:p><w:p><w:r><w:t>• Implement core operations: addition, subtraction, multiplication, division.</w:t></w:r></w:p><w:p><w:r><w:t>• Handle errors such as division by zero and invalid inputs.</w:t></w:r></w:p><w:p><w:r><w:t>• Expand the calculator to include scientific operations (square root, power, etc.).</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• Add a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tkinter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> GUI for a user-friendly experience.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• Integrate </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for version control and documentation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>3. Features (Planned)</w:t></w:r></w:p><w:p><w:r><w:t>• Basic arithmetic operations: +, -, *, /</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• Scientific functions: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sqrt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, power, trigonometric functions</w:t></w:r></w:p><w:p><w:r><w:t>• History of previous calculations</w:t></w:r></w:p><w:p><w:r><w:t>• Error handling for invalid inputs</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tkinter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> GUI with buttons and display area</w:t></w:r></w:p><w:p><w:r><w:t>• Dark mode and light mode (optional)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>4. Project Structure</w:t></w:r></w:p><w:p><w:r><w:br/><w:t xml:space="preserve">    The project will be organized as follows:</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    ├── </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>calculator_project</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    │   ├── main.py                # Entry point for the application</w:t></w:r><w:r><w:br/></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">    │   ├── modules/</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    │   │   ├── __init__.py</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    │   │   ├── mat_utils.py       # Contains mathematical functions</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    │   │   └── sci_utils.py       # Contains scientific functions (future)</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    │   └── </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>gui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>/</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    │       └── interface.py       # </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tkinter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> GUI implementation (future)</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">    </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>5. Development Roadmap</w:t></w:r></w:p><w:p><w:r><w:t>• Phase 1: Implement core operations in mat_utils.py</w:t></w:r></w:p><w:p><w:r><w:t>• Phase 2: Create main.py to call these operations</w:t></w:r></w:p><w:p><w:r><w:t>• Phase 3: Add error handling and unit tests</w:t></w:r></w:p><w:p><w:r><w:t>• Phase 4: Expand with scientific functions (sci_utils.py)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• Phase 5: Build </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tkinter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> GUI for the calculator</w:t></w:r></w:p><w:p><w:r><w:t>• Phase 6: Polish the app with themes and documentation</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$d.Content.InsertXML($xmlFrag)
